$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-6 (columns B-K); column A (index) stays the same.
$data = @{
    2 = @{ B = 1.901710291787398;  C = 4185.8;   D = 0.01349240180076153; E = 29.6; F = 204.4; G = "MyDogaN";     H = "DUO_SUPPORT"; I = 0.09128055152368671; J = 18.6; K = 0.008167217339014521 }
    3 = @{ B = 2.791646791513082;  C = 3890.6;   D = 0.01949175076877424; E = 28.2; F = 221.4; G = "Mr Kayn";     H = "DUO_SUPPORT"; I = 0.1749684198889241;  J = 3.2;  K = 0.002463335941977546 }
    4 = @{ B = 2.416763848396501;  C = 3315.8;   D = 0.0163265306122449;  E = 22.4; F = 40.6;  G = "Booogeyman";  H = "DUO_CARRY";   I = 0.02959183673469388; J = 3.4;  K = 0.002478134110787172 }
    5 = @{ B = 9.48417331160671;   C = 15096.8;  D = 0.04245232838056302; E = 67;   F = 161.6; G = "Shiller";     H = "DUO_CARRY";   I = 0.1032920977440951;  J = 5.4;  K = 0.003379839878639451 }
    6 = @{ B = 1.249264705882353;  C = 2038.8;   D = 0.01409313725490196; E = 23;   F = 138;   G = "Poppy Gods";  H = "DUO_CARRY";   I = 0.08455882352941177; J = 3;    K = 0.001838235294117647 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("B$row").Value = $cols.B
    $ws.Range("C$row").Value = $cols.C
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("F$row").Value = $cols.F
    $ws.Range("G$row").Value = $cols.G
    $ws.Range("H$row").Value = $cols.H
    $ws.Range("I$row").Value = $cols.I
    $ws.Range("J$row").Value = $cols.J
    $ws.Range("K$row").Value = $cols.K
}
